$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 176  # H4: 156 -> 176
$ws.Cells.Item(4, 9).Value = 170  # I4: 156 -> 170
$ws.Cells.Item(4, 10).Value = 200  # J4: 0 -> 200
$ws.Cells.Item(4, 11).Value = 170  # K4: 156 -> 170
$ws.Cells.Item(4, 12).Value = 200  # L4: 0 -> 200
$ws.Cells.Item(4, 13).Value = -56  # M4: -42 -> -56
$ws.Cells.Item(4, 14).Value = -428  # N4: None -> -428
$ws.Cells.Item(19, 8).Value = 1224.2858  # H19: 1210.8462 -> 1224.2858
$ws.Cells.Item(19, 9).Value = 1009.7  # I19: 966.3333 -> 1009.7
$ws.Cells.Item(19, 10).Value = 1760.75  # J19: 1761 -> 1760.75
$ws.Cells.Item(19, 11).Value = 1009.7  # K19: 966.3333 -> 1009.7
$ws.Cells.Item(19, 12).Value = 1760.75  # L19: 1761 -> 1760.75
$ws.Cells.Item(19, 13).Value = -834.7  # M19: -791.3333 -> -834.7
$ws.Cells.Item(19, 14).Value = -2110.75  # N19: -2111 -> -2110.75
$ws.Cells.Item(43, 8).Value = 3333  # H43: 3499.5 -> 3333
$ws.Cells.Item(43, 10).Value = 3000  # J43: 0 -> 3000
$ws.Cells.Item(43, 12).Value = 3000  # L43: 0 -> 3000
$ws.Cells.Item(43, 14).Value = -3138  # N43: None -> -3138
$ws.Cells.Item(58, 8).Value = 329.8  # H58: 289.83334 -> 329.8
$ws.Cells.Item(58, 9).Value = 329.8  # I58: 289.83334 -> 329.8
$ws.Cells.Item(58, 11).Value = 989.4000000000001  # K58: 869.5000200000001 -> 989.4000000000001
$ws.Cells.Item(58, 13).Value = -839.4000000000001  # M58: -719.5000200000001 -> -839.4000000000001
$ws.Cells.Item(64, 8).Value = 4981.3335  # H64: 5000 -> 4981.3335
$ws.Cells.Item(64, 9).Value = 4977.6  # I64: 5000 -> 4977.6
$ws.Cells.Item(64, 11).Value = 4977.6  # K64: 5000 -> 4977.6
$ws.Cells.Item(64, 13).Value = -4729.6  # M64: -4752 -> -4729.6
$ws.Cells.Item(67, 8).Value = 4981.3335  # H67: 5000 -> 4981.3335
$ws.Cells.Item(67, 9).Value = 4977.6  # I67: 5000 -> 4977.6
$ws.Cells.Item(67, 11).Value = 4977.6  # K67: 5000 -> 4977.6
$ws.Cells.Item(67, 13).Value = -4119.6  # M67: -4142 -> -4119.6
$ws.Cells.Item(116, 8).Value = 4392  # H116: 3718.75 -> 4392
$ws.Cells.Item(116, 9).Value = 3632.6667  # I116: 2859.2 -> 3632.6667
$ws.Cells.Item(116, 11).Value = 3632.6667  # K116: 2859.2 -> 3632.6667
$ws.Cells.Item(116, 13).Value = -190.6667000000002  # M116: 582.8000000000002 -> -190.6667000000002
$ws.Cells.Item(131, 8).Value = 300  # H131: 1649.5 -> 300
$ws.Cells.Item(131, 10).Value = 0  # J131: 2999 -> 0
$ws.Cells.Item(131, 12).Value = 0  # L131: 8997 -> 0
$ws.Cells.Item(131, 14).ClearContents()  # N131 was -19077

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 274.75  # H5: 269.8 -> 274.75
$ws.Cells.Item(5, 9).Value = 266.66666  # I5: 250 -> 266.66666
$ws.Cells.Item(5, 10).Value = 299  # J5: 299.5 -> 299
$ws.Cells.Item(5, 11).Value = 266.66666  # K5: 250 -> 266.66666
$ws.Cells.Item(5, 12).Value = 299  # L5: 299.5 -> 299
$ws.Cells.Item(5, 13).Value = -154.66666  # M5: -138 -> -154.66666
$ws.Cells.Item(5, 14).Value = -523  # N5: -523.5 -> -523
$ws.Cells.Item(45, 8).Value = 1999  # H45: 1887.5 -> 1999
$ws.Cells.Item(45, 9).Value = 0  # I45: 1775 -> 0
$ws.Cells.Item(45, 10).Value = 1999  # J45: 2000 -> 1999
$ws.Cells.Item(45, 11).Value = 0  # K45: 1775 -> 0
$ws.Cells.Item(45, 12).Value = 1999  # L45: 2000 -> 1999
$ws.Cells.Item(45, 13).ClearContents()  # M45 was -1398
$ws.Cells.Item(45, 14).Value = -2753  # N45: -2754 -> -2753
$ws.Cells.Item(61, 8).Value = 7333  # H61: 8000 -> 7333
$ws.Cells.Item(61, 9).Value = 7333  # I61: 8000 -> 7333
$ws.Cells.Item(61, 11).Value = 7333  # K61: 8000 -> 7333
$ws.Cells.Item(61, 13).Value = -7121  # M61: -7788 -> -7121
$ws.Cells.Item(105, 8).Value = 37500  # H105: 0 -> 37500
$ws.Cells.Item(105, 10).Value = 37500  # J105: 0 -> 37500
$ws.Cells.Item(105, 12).Value = 37500  # L105: 0 -> 37500
$ws.Cells.Item(105, 14).Value = -44488  # N105: None -> -44488
$ws.Cells.Item(132, 8).Value = 1029.6666  # H132: 902.9 -> 1029.6666
$ws.Cells.Item(132, 9).Value = 1029.6666  # I132: 836.55554 -> 1029.6666
$ws.Cells.Item(132, 10).Value = 0  # J132: 1500 -> 0
$ws.Cells.Item(132, 11).Value = 3088.9998  # K132: 2509.66662 -> 3088.9998
$ws.Cells.Item(132, 12).Value = 0  # L132: 4500 -> 0
$ws.Cells.Item(132, 13).Value = -558.9998000000001  # M132: 20.33338000000003 -> -558.9998000000001
$ws.Cells.Item(132, 14).ClearContents()  # N132 was -9560
$ws.Cells.Item(136, 8).Value = 7333  # H136: 8000 -> 7333
$ws.Cells.Item(136, 9).Value = 7333  # I136: 8000 -> 7333
$ws.Cells.Item(136, 11).Value = 21999  # K136: 24000 -> 21999
$ws.Cells.Item(136, 13).Value = -19449  # M136: -21450 -> -19449

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 274.75  # H4: 269.8 -> 274.75
$ws.Cells.Item(4, 9).Value = 266.66666  # I4: 250 -> 266.66666
$ws.Cells.Item(4, 10).Value = 299  # J4: 299.5 -> 299
$ws.Cells.Item(4, 11).Value = 266.66666  # K4: 250 -> 266.66666
$ws.Cells.Item(4, 12).Value = 299  # L4: 299.5 -> 299
$ws.Cells.Item(4, 13).Value = -151.66666  # M4: -135 -> -151.66666
$ws.Cells.Item(4, 14).Value = -529  # N4: -529.5 -> -529
$ws.Cells.Item(86, 8).Value = 4466.5557  # H86: 3339.8 -> 4466.5557
$ws.Cells.Item(86, 9).Value = 2538.4  # I86: 2923 -> 2538.4
$ws.Cells.Item(86, 10).Value = 6876.75  # J86: 5007 -> 6876.75
$ws.Cells.Item(86, 11).Value = 2538.4  # K86: 2923 -> 2538.4
$ws.Cells.Item(86, 12).Value = 6876.75  # L86: 5007 -> 6876.75
$ws.Cells.Item(86, 13).Value = -1415.4  # M86: -1800 -> -1415.4
$ws.Cells.Item(86, 14).Value = -9122.75  # N86: -7253 -> -9122.75
$ws.Cells.Item(89, 8).Value = 4466.5557  # H89: 3339.8 -> 4466.5557
$ws.Cells.Item(89, 9).Value = 2538.4  # I89: 2923 -> 2538.4
$ws.Cells.Item(89, 10).Value = 6876.75  # J89: 5007 -> 6876.75
$ws.Cells.Item(89, 11).Value = 12692  # K89: 14615 -> 12692
$ws.Cells.Item(89, 12).Value = 34383.75  # L89: 25035 -> 34383.75
$ws.Cells.Item(89, 13).Value = -7076  # M89: -8999 -> -7076
$ws.Cells.Item(89, 14).Value = -45615.75  # N89: -36267 -> -45615.75
$ws.Cells.Item(99, 8).Value = 6560  # H99: 7600 -> 6560
$ws.Cells.Item(99, 9).Value = 6560  # I99: 7600 -> 6560
$ws.Cells.Item(99, 11).Value = 6560  # K99: 7600 -> 6560
$ws.Cells.Item(99, 13).Value = -5062  # M99: -6102 -> -5062
$ws.Cells.Item(134, 8).Value = 3950  # H134: 4573.1177 -> 3950
$ws.Cells.Item(134, 9).Value = 3647.1177  # I134: 4338.857 -> 3647.1177
$ws.Cells.Item(134, 11).Value = 10941.3531  # K134: 13016.571 -> 10941.3531
$ws.Cells.Item(134, 13).Value = -8406.3531  # M134: -10481.571 -> -8406.3531

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 800  # H4: 0 -> 800
$ws.Cells.Item(4, 9).Value = 800  # I4: 0 -> 800
$ws.Cells.Item(4, 11).Value = 800  # K4: 0 -> 800
$ws.Cells.Item(4, 13).Value = -688  # M4: None -> -688
$ws.Cells.Item(45, 8).Value = 52000  # H45: 0 -> 52000
$ws.Cells.Item(45, 10).Value = 52000  # J45: 0 -> 52000
$ws.Cells.Item(45, 12).Value = 52000  # L45: 0 -> 52000
$ws.Cells.Item(45, 14).Value = -53186  # N45: None -> -53186
$ws.Cells.Item(58, 8).Value = 3598  # H58: 2528.3333 -> 3598
$ws.Cells.Item(58, 9).Value = 3598  # I58: 2528.3333 -> 3598
$ws.Cells.Item(58, 11).Value = 3598  # K58: 2528.3333 -> 3598
$ws.Cells.Item(58, 13).Value = -3395  # M58: -2325.3333 -> -3395
$ws.Cells.Item(122, 8).Value = 1393.6666  # H122: 1376.8 -> 1393.6666
$ws.Cells.Item(122, 9).Value = 1450.4286  # I122: 1422.25 -> 1450.4286
$ws.Cells.Item(122, 11).Value = 4351.2858  # K122: 4266.75 -> 4351.2858
$ws.Cells.Item(122, 13).Value = -1901.2858  # M122: -1816.75 -> -1901.2858
$ws.Cells.Item(134, 8).Value = 2000  # H134: 0 -> 2000
$ws.Cells.Item(134, 10).Value = 2000  # J134: 0 -> 2000
$ws.Cells.Item(134, 12).Value = 6000  # L134: 0 -> 6000
$ws.Cells.Item(134, 14).Value = -11070  # N134: None -> -11070
$ws.Cells.Item(136, 8).Value = 3598  # H136: 2528.3333 -> 3598
$ws.Cells.Item(136, 9).Value = 3598  # I136: 2528.3333 -> 3598
$ws.Cells.Item(136, 11).Value = 10794  # K136: 7584.999899999999 -> 10794
$ws.Cells.Item(136, 13).Value = -8244  # M136: -5034.999899999999 -> -8244

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 97431.67999999999  # H4: 89143.25999999999 -> 97431.67999999999
$ws.Cells.Item(4, 9).Value = 176733.33  # I4: 168750 -> 176733.33
$ws.Cells.Item(4, 10).Value = 2269.7  # J4: 2299.5454 -> 2269.7
$ws.Cells.Item(4, 11).Value = 530199.99  # K4: 506250 -> 530199.99
$ws.Cells.Item(4, 12).Value = 6809.099999999999  # L4: 6898.6362 -> 6809.099999999999
$ws.Cells.Item(4, 13).Value = -530087.99  # M4: -506138 -> -530087.99
$ws.Cells.Item(4, 14).Value = -7033.099999999999  # N4: -7122.6362 -> -7033.099999999999
$ws.Cells.Item(9, 8).Value = 333  # H9: 381.125 -> 333
$ws.Cells.Item(9, 9).Value = 392.42856  # I9: 441.66666 -> 392.42856
$ws.Cells.Item(9, 10).Value = 125  # J9: 344.8 -> 125
$ws.Cells.Item(9, 11).Value = 1177.28568  # K9: 1324.99998 -> 1177.28568
$ws.Cells.Item(9, 12).Value = 375  # L9: 1034.4 -> 375
$ws.Cells.Item(9, 13).Value = -953.28568  # M9: -1100.99998 -> -953.28568
$ws.Cells.Item(9, 14).Value = -823  # N9: -1482.4 -> -823
$ws.Cells.Item(68, 8).Value = 2049.75  # H68: 2200 -> 2049.75
$ws.Cells.Item(68, 9).Value = 1799.5  # I68: 2000 -> 1799.5
$ws.Cells.Item(68, 11).Value = 5398.5  # K68: 6000 -> 5398.5
$ws.Cells.Item(68, 13).Value = -4587.5  # M68: -5189 -> -4587.5
$ws.Cells.Item(71, 8).Value = 2049.75  # H71: 2200 -> 2049.75
$ws.Cells.Item(71, 9).Value = 1799.5  # I71: 2000 -> 1799.5
$ws.Cells.Item(71, 11).Value = 16195.5  # K71: 18000 -> 16195.5
$ws.Cells.Item(71, 13).Value = -12139.5  # M71: -13944 -> -12139.5
$ws.Cells.Item(92, 8).Value = 999.4286  # H92: 1250 -> 999.4286
$ws.Cells.Item(92, 9).Value = 251  # I92: 0 -> 251
$ws.Cells.Item(92, 10).Value = 1298.8  # J92: 1250 -> 1298.8
$ws.Cells.Item(92, 11).Value = 753  # K92: 0 -> 753
$ws.Cells.Item(92, 12).Value = 3896.4  # L92: 3750 -> 3896.4
$ws.Cells.Item(92, 13).Value = 495  # M92: None -> 495
$ws.Cells.Item(92, 14).Value = -6392.4  # N92: -6246 -> -6392.4
$ws.Cells.Item(98, 8).Value = 297  # H98: 0 -> 297
$ws.Cells.Item(98, 10).Value = 297  # J98: 0 -> 297
$ws.Cells.Item(98, 12).Value = 891  # L98: 0 -> 891
$ws.Cells.Item(98, 14).Value = -3887  # N98: None -> -3887

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 360  # H2: 103.1875 -> 360
$ws.Cells.Item(2, 9).Value = 53.6  # I2: 99.8 -> 53.6
$ws.Cells.Item(2, 10).Value = 797.7143  # J2: 108.833336 -> 797.7143
$ws.Cells.Item(2, 11).Value = 53.6  # K2: 99.8 -> 53.6
$ws.Cells.Item(2, 12).Value = 797.7143  # L2: 108.833336 -> 797.7143
$ws.Cells.Item(2, 13).Value = 59.4  # M2: 13.2 -> 59.4
$ws.Cells.Item(2, 14).Value = -1023.7143  # N2: -334.833336 -> -1023.7143
$ws.Cells.Item(11, 8).Value = 133500  # H11: 3400166.8 -> 133500
$ws.Cells.Item(11, 9).Value = 200000  # I11: 0 -> 200000
$ws.Cells.Item(11, 10).Value = 100250  # J11: 3400166.8 -> 100250
$ws.Cells.Item(11, 11).Value = 200000  # K11: 0 -> 200000
$ws.Cells.Item(11, 12).Value = 100250  # L11: 3400166.8 -> 100250
$ws.Cells.Item(11, 13).Value = -199861  # M11: None -> -199861
$ws.Cells.Item(11, 14).Value = -100528  # N11: -3400444.8 -> -100528
$ws.Cells.Item(102, 8).Value = 1493.2222  # H102: 1723 -> 1493.2222
$ws.Cells.Item(102, 10).Value = 1514  # J102: 2064 -> 1514
$ws.Cells.Item(102, 12).Value = 1514  # L102: 2064 -> 1514
$ws.Cells.Item(102, 14).Value = -4758  # N102: -5308 -> -4758
$ws.Cells.Item(122, 8).Value = 0  # H122: 500 -> 0
$ws.Cells.Item(122, 9).Value = 0  # I122: 500 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 1500 -> 0
$ws.Cells.Item(122, 13).ClearContents()  # M122 was 950
$ws.Cells.Item(126, 8).Value = 8698.817999999999  # H126: 8934.727999999999 -> 8698.817999999999
$ws.Cells.Item(126, 9).Value = 8698.817999999999  # I126: 8934.727999999999 -> 8698.817999999999
$ws.Cells.Item(126, 11).Value = 26096.454  # K126: 26804.184 -> 26096.454
$ws.Cells.Item(126, 13).Value = -23626.454  # M126: -24334.184 -> -23626.454
$ws.Cells.Item(132, 8).Value = 1491.5  # H132: 1494.3334 -> 1491.5
$ws.Cells.Item(132, 9).Value = 1491.5  # I132: 1494.3334 -> 1491.5
$ws.Cells.Item(132, 11).Value = 4474.5  # K132: 4483.0002 -> 4474.5
$ws.Cells.Item(132, 13).Value = -1944.5  # M132: -1953.0002 -> -1944.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 23230.578  # H7: 25576.5 -> 23230.578
$ws.Cells.Item(7, 9).Value = 24213.646  # I7: 26914.25 -> 24213.646
$ws.Cells.Item(7, 11).Value = 24213.646  # K7: 26914.25 -> 24213.646
$ws.Cells.Item(7, 13).Value = -24101.646  # M7: -26802.25 -> -24101.646
$ws.Cells.Item(22, 8).Value = 1949.6666  # H22: 2000 -> 1949.6666
$ws.Cells.Item(22, 9).Value = 1850  # I22: 0 -> 1850
$ws.Cells.Item(22, 10).Value = 1999.5  # J22: 2000 -> 1999.5
$ws.Cells.Item(22, 11).Value = 1850  # K22: 0 -> 1850
$ws.Cells.Item(22, 12).Value = 1999.5  # L22: 2000 -> 1999.5
$ws.Cells.Item(22, 13).Value = -1555  # M22: None -> -1555
$ws.Cells.Item(22, 14).Value = -2589.5  # N22: -2590 -> -2589.5
$ws.Cells.Item(27, 8).Value = 1949.6666  # H27: 2000 -> 1949.6666
$ws.Cells.Item(27, 9).Value = 1850  # I27: 0 -> 1850
$ws.Cells.Item(27, 10).Value = 1999.5  # J27: 2000 -> 1999.5
$ws.Cells.Item(27, 11).Value = 1850  # K27: 0 -> 1850
$ws.Cells.Item(27, 12).Value = 1999.5  # L27: 2000 -> 1999.5
$ws.Cells.Item(27, 13).Value = -1743  # M27: None -> -1743
$ws.Cells.Item(27, 14).Value = -2213.5  # N27: -2214 -> -2213.5
$ws.Cells.Item(68, 8).Value = 28911.25  # H68: 32185 -> 28911.25
$ws.Cells.Item(68, 10).Value = 54498.25  # J68: 70666 -> 54498.25
$ws.Cells.Item(68, 12).Value = 54498.25  # L68: 70666 -> 54498.25
$ws.Cells.Item(68, 14).Value = -55996.25  # N68: -72164 -> -55996.25
$ws.Cells.Item(71, 8).Value = 28911.25  # H71: 32185 -> 28911.25
$ws.Cells.Item(71, 10).Value = 54498.25  # J71: 70666 -> 54498.25
$ws.Cells.Item(71, 12).Value = 272491.25  # L71: 353330 -> 272491.25
$ws.Cells.Item(71, 14).Value = -279979.25  # N71: -360818 -> -279979.25
$ws.Cells.Item(125, 8).Value = 0  # H125: 15000 -> 0
$ws.Cells.Item(125, 10).Value = 0  # J125: 15000 -> 0
$ws.Cells.Item(125, 12).Value = 0  # L125: 15000 -> 0
$ws.Cells.Item(125, 14).ClearContents()  # N125 was -24840
$ws.Cells.Item(126, 8).Value = 23230.578  # H126: 25576.5 -> 23230.578
$ws.Cells.Item(126, 9).Value = 24213.646  # I126: 26914.25 -> 24213.646
$ws.Cells.Item(126, 11).Value = 72640.93799999999  # K126: 80742.75 -> 72640.93799999999
$ws.Cells.Item(126, 13).Value = -70170.93799999999  # M126: -78272.75 -> -70170.93799999999
$ws.Cells.Item(136, 8).Value = 2169.8333  # H136: 2084.8333 -> 2169.8333
$ws.Cells.Item(136, 9).Value = 2003.3529  # I136: 2084.8333 -> 2003.3529
$ws.Cells.Item(136, 10).Value = 5000  # J136: 0 -> 5000
$ws.Cells.Item(136, 11).Value = 6010.0587  # K136: 6254.499899999999 -> 6010.0587
$ws.Cells.Item(136, 12).Value = 15000  # L136: 0 -> 15000
$ws.Cells.Item(136, 13).Value = -3460.0587  # M136: -3704.499899999999 -> -3460.0587
$ws.Cells.Item(136, 14).Value = -20100  # N136: None -> -20100

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 9997  # H2: 12666.333 -> 9997
$ws.Cells.Item(2, 9).Value = 0  # I2: 8002 -> 0
$ws.Cells.Item(2, 10).Value = 9997  # J2: 14998.5 -> 9997
$ws.Cells.Item(2, 11).Value = 0  # K2: 8002 -> 0
$ws.Cells.Item(2, 12).Value = 9997  # L2: 14998.5 -> 9997
$ws.Cells.Item(2, 13).ClearContents()  # M2 was -7890
$ws.Cells.Item(2, 14).Value = -10221  # N2: -15222.5 -> -10221
$ws.Cells.Item(132, 8).Value = 2427.7693  # H132: 2560.0908 -> 2427.7693
$ws.Cells.Item(132, 9).Value = 2421.8333  # I132: 2566.2 -> 2421.8333
$ws.Cells.Item(132, 11).Value = 7265.499899999999  # K132: 7698.599999999999 -> 7265.499899999999
$ws.Cells.Item(132, 13).Value = -4735.499899999999  # M132: -5168.599999999999 -> -4735.499899999999
